$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the updated "Price" values are plain decimal numbers (e.g. "300.96").
# Assigning such a string straight to .Value lets Excel auto-convert it to a
# real number, which both changes the cell type and introduces binary
# floating-point rounding noise (e.g. 97.88 -> 97.879999999999995). The source
# data must stay plain text, exactly as authored, so those cells are switched to
# Text number format before the value is written, then restored to the default
# (Normal) style afterwards so no stray formatting is left behind.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D13", "D14", "D19", "D21", "D23", "D25", "D26", "D28", "D31", "D33", "D34", "D35", "D37", "D38", "D44", "D45", "D46", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.057.59'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '2.303.23'
$ws.Range("E3").Value = '  +0.45%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '300.96'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").Value = '97.88'
$ws.Range("E6").Value = '  -1.16%  '
$ws.Range("D7").Value = '0.522'
$ws.Range("E7").Value = '  +4.36%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").Value = '0.517'
$ws.Range("E9").Value = '  +0.83%  '
$ws.Range("D10").Value = '35.73'
$ws.Range("E10").Value = '  -0.96%  '
$ws.Range("D11").Value = '0.0790'
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("D13").Value = '18.00'
$ws.Range("E13").Value = '  +0.76%  '
$ws.Range("D14").Value = '6.86'
$ws.Range("E14").Value = '  +0.24%  '
$ws.Range("D15").Value = '2.662.59'
$ws.Range("E15").Value = '  +0.54%  '
$ws.Range("D16").Value = '2.314.11'
$ws.Range("E16").Value = '  -0.78%  '
$ws.Range("E17").Value = '  -2.04%  '
$ws.Range("D18").Value = '42.966.22'
$ws.Range("E18").Value = '  +0.67%  '
$ws.Range("D19").Value = '13.37'
$ws.Range("E19").Value = '  +8.33%  '
$ws.Range("E20").Value = '  +0.86%  '
$ws.Range("D21").Value = '6.11'
$ws.Range("E21").Value = '  -1.27%  '
$ws.Range("E22").Value = '  +0.64%  '
$ws.Range("D23").Value = '238.31'
$ws.Range("E23").Value = '  +1.01%  '
$ws.Range("E24").Value = '  -1.49%  '
$ws.Range("D25").Value = '0.998'
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("D26").Value = '2.43'
$ws.Range("E26").Value = '  -0.50%  '
$ws.Range("E27").Value = '  +0.36%  '
$ws.Range("D28").Value = '167.83'
$ws.Range("E28").Value = '  -0.46%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  -0.90%  '
$ws.Range("D31").Value = '32.92'
$ws.Range("E31").Value = '  -4.60%  '
$ws.Range("E32").Value = '  +4.01%  '
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").Value = '4.80'
$ws.Range("E34").Value = '  +4.67%  '
$ws.Range("D35").Value = '18.07'
$ws.Range("E35").Value = '  +3.47%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("D37").Value = '0.0688'
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("D38").Value = '0.102'
$ws.Range("E39").Value = '  +0.85%  '
$ws.Range("E40").Value = '  +2.21%  '
$ws.Range("E41").Value = '  -2.65%  '
$ws.Range("D42").Value = '2.013.24'
$ws.Range("E42").Value = '  +0.73%  '
$ws.Range("E43").Value = '  +0.66%  '
$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").Value = '2.16'
$ws.Range("E44").Value = '  -2.86%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '10.19'
$ws.Range("E45").Value = '  +0.63%  '
$ws.Range("D46").Value = '17.32'
$ws.Range("E46").Value = '  -0.81%  '
$ws.Range("E47").Value = '  -1.91%  '
$ws.Range("D48").Value = '54.32'
$ws.Range("E48").Value = '  -2.21%  '
$ws.Range("D49").Value = '2.528.90'
$ws.Range("E49").Value = '  +0.59%  '
$ws.Range("E50").Value = '  +0.68%  '
$ws.Range("D51").Value = '73.19'
$ws.Range("E51").Value = '  +4.38%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
